# Update the LR-pair worksheet with new TPM-derived values.
# The "ECs" sending-cluster rows are removed entirely, the remaining
# FAPs/MuSCs sending-cluster rows move up (now rows 2-4 and 5-7), and
# all numeric columns (E:T) are recalculated with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data: 6 rows (FAPs->{ECs,FAPs,MuSCs}, MuSCs->{ECs,FAPs,MuSCs}) ---
$data = @(
    @("FAPs",  "Alcam", "L1cam", "ECs",   3, 1, 0.5683613333333334, 1.705084, 0.4361027177196302, 0.4361027177196302, 3, 1, 3.685507, 11.056521, 0.3585631737883472, 0.3585631737883472, 2.094699672529333,  18.852297052764,  0.1563703745632743, 0.1563703745632743),
    @("FAPs",  "Alcam", "L1cam", "FAPs",  3, 1, 0.5683613333333334, 1.705084, 0.4361027177196302, 0.4361027177196302, 1, 0.3333333333333333, 0.099159, 0.297477, 0.009647184430711629, 0.009647184430711629, 0.056358141452, 0.507223273068, 0.004207163348575845, 0.004207163348575845),
    @("FAPs",  "Alcam", "L1cam", "MuSCs", 3, 1, 0.5683613333333334, 1.705084, 0.4361027177196302, 0.4361027177196302, 3, 1, 6.493877, 19.481631, 0.6317896417809412, 0.6317896417809411, 3.690868590222667,  33.217817312004,  0.2755251798077801, 0.2755251798077801),
    @("MuSCs", "Alcam", "L1cam", "ECs",   3, 1, 0.7349126666666667, 2.204738, 0.5638972822803697, 0.5638972822803697, 3, 1, 3.685507, 11.056521, 0.3585631737883472, 0.3585631737883472, 2.708525777388667,  24.376731996498,  0.2021927992250729, 0.2021927992250729),
    @("MuSCs", "Alcam", "L1cam", "FAPs",  3, 1, 0.7349126666666667, 2.204738, 0.5638972822803697, 0.5638972822803697, 1, 0.3333333333333333, 0.099159, 0.297477, 0.009647184430711629, 0.009647184430711629, 0.072873205114, 0.655858846026, 0.005440021082135783, 0.005440021082135783),
    @("MuSCs", "Alcam", "L1cam", "MuSCs", 3, 1, 0.7349126666666667, 2.204738, 0.5638972822803697, 0.5638972822803697, 3, 1, 6.493877, 19.481631, 0.6317896417809412, 0.6317896417809411, 4.772432463075334,  42.951892167678,  0.3562644619731611, 0.356264461973161)
)

# Remove the now-obsolete rows 8, 9, 10 (old MuSCs-as-sender block that
# trails the data after the ECs block is dropped and the rest shift up).
$ws.Rows.Item(8).Resize(3).EntireRow.Delete()

# Write the 6 rows of new data into A2:T7, column by column (matches the
# source script's generation order, which iterates per-column).
$numCols = $data[0].Length
for ($c = 0; $c -lt $numCols; $c++) {
    for ($i = 0; $i -lt $data.Length; $i++) {
        $r = $i + 2
        $ws.Cells.Item($r, $c + 1).Value = $data[$i][$c]
    }
}
